$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (B & C) for the new "Jun_17" / "Jun_15" report
#    dates. This shifts the old "Jun_13" column (B) to D and the old
#    "Jun_10" column (C, which also carried the per-analyst rating-change
#    notes) to E -- exactly mirroring what MarketBeat's scraper does each
#    time it appends a newer snapshot to the left of the table.
# ---------------------------------------------------------------------------
$ws.Range("B:C").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. Header row: newest dates go in the freshly inserted columns.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# ---------------------------------------------------------------------------
# 3. Fill the new columns with the default "UN" (unchanged) rating for
#    every analyst row.
# ---------------------------------------------------------------------------
$ws.Range("B2:C27").Value = "UN"

# ---------------------------------------------------------------------------
# 4. Record the one new rating action that happened in this snapshot: on
#    6/15/2018 the rating on row 11 (ValuEngine) was downgraded.
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = "6/15/2018,Downgrades,Hold -> Sell,"
$ws.Range("C11").Interior.ColorIndex = 45
$ws.Range("C11").Interior.Pattern = -4142
$ws.Range("C11").Interior.Pattern = 1

# ---------------------------------------------------------------------------
# 5. Column widths: keep the historical "Jun_10" column (now E) at its
#    explicit width, and carry the same width onto the two columns that
#    used to sit to its left (C, D) before they got pushed further along.
# ---------------------------------------------------------------------------
$ws.Range("C:C").ColumnWidth = 8
$ws.Range("D:D").ColumnWidth = 8
$ws.Range("E:E").ColumnWidth = 8
